$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (rows 2..30), describing how the
# weekly data set was reshuffled across the date-ordered rows. The full set
# of 29 data rows is identical before/after; only their row positions change.
$rowMap = @{
    2 = 17
    3 = 18
    4 = 19
    5 = 20
    6 = 6
    7 = 28
    8 = 29
    9 = 30
    10 = 14
    11 = 15
    12 = 2
    13 = 3
    14 = 26
    15 = 27
    16 = 21
    17 = 22
    18 = 10
    19 = 8
    20 = 9
    21 = 13
    22 = 4
    23 = 5
    24 = 12
    25 = 24
    26 = 25
    27 = 11
    28 = 7
    29 = 16
    30 = 23
}

# Snapshot all source values (columns A..T, rows 2..30) before writing,
# since this is a permutation and rows both read from and written to.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    for ($c = 1; $c -le 20; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot["$srcRow,$c"]
    }
}
